$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 33.84781966666667
$ws.Range("H2").Value = 101.543459
$ws.Range("I2").Value = 0.2402182618707165
$ws.Range("J2").Value = 0.2402182618707166
$ws.Range("M2").Value = 0.4890553333333333
$ws.Range("N2").Value = 1.467166
$ws.Range("O2").Value = 0.9644476581758422
$ws.Range("P2").Value = 0.9644476581758422
$ws.Range("Q2").Value = 16.55345672968822
$ws.Range("R2").Value = 148.981110567194
$ws.Range("S2").Value = 0.2316779401122838
$ws.Range("T2").Value = 0.2316779401122838

# Row 3
$ws.Range("G3").Value = 33.84781966666667
$ws.Range("H3").Value = 101.543459
$ws.Range("I3").Value = 0.2402182618707165
$ws.Range("J3").Value = 0.2402182618707166
$ws.Range("O3").Value = 0.03555234182415776
$ws.Range("P3").Value = 0.03555234182415776
$ws.Range("Q3").Value = 0.6102084929506666
$ws.Range("R3").Value = 5.491876436556
$ws.Range("S3").Value = 0.008540321758432758
$ws.Range("T3").Value = 0.008540321758432759

# Row 4
$ws.Range("I4").Value = 0.5284143281787288
$ws.Range("J4").Value = 0.5284143281787288
$ws.Range("M4").Value = 0.4890553333333333
$ws.Range("N4").Value = 1.467166
$ws.Range("O4").Value = 0.9644476581758422
$ws.Range("P4").Value = 0.9644476581758422
$ws.Range("Q4").Value = 36.41306721951666
$ws.Range("R4").Value = 327.71760497565
$ws.Range("S4").Value = 0.5096279613585359
$ws.Range("T4").Value = 0.5096279613585359

# Row 5
$ws.Range("I5").Value = 0.5284143281787288
$ws.Range("J5").Value = 0.5284143281787288
$ws.Range("O5").Value = 0.03555234182415776
$ws.Range("P5").Value = 0.03555234182415776
$ws.Range("S5").Value = 0.01878636682019285
$ws.Range("T5").Value = 0.01878636682019285

# Row 6
$ws.Range("H6").Value = 97.802086
$ws.Range("I6").Value = 0.2313674099505547
$ws.Range("J6").Value = 0.2313674099505547
$ws.Range("M6").Value = 0.4890553333333333
$ws.Range("N6").Value = 1.467166
$ws.Range("O6").Value = 0.9644476581758422
$ws.Range("P6").Value = 0.9644476581758422
$ws.Range("Q6").Value = 15.94354392314178
$ws.Range("R6").Value = 143.491895308276
$ws.Range("S6").Value = 0.2231417567050225
$ws.Range("T6").Value = 0.2231417567050225

# Row 7
$ws.Range("H7").Value = 97.802086
$ws.Range("I7").Value = 0.2313674099505547
$ws.Range("J7").Value = 0.2313674099505547
$ws.Range("O7").Value = 0.03555234182415776
$ws.Range("P7").Value = 0.03555234182415776
$ws.Range("Q7").Value = 0.5877253354693334
$ws.Range("R7").Value = 5.289528019224
$ws.Range("S7").Value = 0.00822565324553216
$ws.Range("T7").Value = 0.008225653245532162
